$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.247.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.494.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.90%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.415"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.01"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.491.13"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.87"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +11.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.202"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "97.282.03"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.147.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.95%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.72"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.495.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.41"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +11.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +17.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.500"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "522.27"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.32"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.74"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "93.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.50"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.690.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.22"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +14.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.79"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +15.39%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.140"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.187"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +13.10%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.586"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.37%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.87"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.76%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.153"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "512.69"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.915"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +12.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.31"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0422"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.63"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.58"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +13.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.47"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.73%  "
